$d = $word.ActiveDocument

$pairs = @(
    @("54×38=", "34×34="),
    @("20×21=", "95×68="),
    @("66×15=", "98×64="),
    @("81×87=", "97×24="),
    @("73×31=", "60×43="),
    @("11×26=", "66×74="),
    @("82×66=", "59×25="),
    @("36×91=", "86×24="),
    @("96×93=", "75×83="),
    @("43×77=", "68×68="),
    @("60×90=", "45×32="),
    @("74×26=", "99×94="),
    @("19×71=", "47×13="),
    @("61×86=", "22×17="),
    @("45×25=", "49×95="),
    @("78×58=", "70×99="),
    @("78×74=", "97×58="),
    @("76×57=", "39×12="),
    @("25×14=", "38×62="),
    @("30×43=", "85×36="),
    @("84×25=", "65×81="),
    @("47×48=", "65×44="),
    @("48×68=", "61×24="),
    @("60×78=", "60×55="),
    @("43×70=", "41×56=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
